# Add a new data row (row 3) to the "Artfynd" worksheet, mirroring the
# structure of the existing row 2 record for Valeriana dioica.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = 3

$ws.Cells.Item($r, 1).Value  = 112276241                                     # A  Id
$ws.Cells.Item($r, 2).Value  = 108492                                        # B  Taxonsorteringsordning
$ws.Cells.Item($r, 3).Value  = "Ovaliderad"                                  # C  Valideringsstatus
$ws.Cells.Item($r, 4).Value  = "VU"                                          # D  Rödlistade
$ws.Cells.Item($r, 5).Value  = 221049                                        # E  TaxonId
$ws.Cells.Item($r, 6).Value  = "Småvänderot"                                 # F  Artnamn
$ws.Cells.Item($r, 7).Value  = "Valeriana dioica"                            # G  Vetenskapligt namn
$ws.Cells.Item($r, 8).Value  = "L."                                          # H  Auktor

# I, J, K, L, N: present in the source as empty text cells.
$ws.Cells.Item($r, 9).Value  = ""                                            # I  Antal
$ws.Cells.Item($r, 10).Value = ""                                            # J  Enhet
$ws.Cells.Item($r, 11).Value = ""                                            # K  Ålder-Stadium
$ws.Cells.Item($r, 12).Value = ""                                            # L  Kön
$ws.Cells.Item($r, 14).Value = ""                                            # N  Metod

$ws.Cells.Item($r, 16).Value = "Köle vägkorsning, 400 m SSO , Sk"            # P  Lokalnamn
$ws.Cells.Item($r, 17).Value = 448505                                        # Q  Ost
$ws.Cells.Item($r, 18).Value = 6185264                                       # R  Nord
$ws.Cells.Item($r, 19).Value = 50                                            # S  Noggrannhet
$ws.Cells.Item($r, 20).Value = "Skåne"                                       # T  Län
$ws.Cells.Item($r, 21).Value = "Kristianstad"                                # U  Kommun
$ws.Cells.Item($r, 22).Value = "Skåne"                                       # V  Provins
$ws.Cells.Item($r, 23).Value = "Maglehem"                                    # W  Församling
$ws.Cells.Item($r, 24).Value = "M-Kri-0111"                                  # X  Externid

# Y, AA hold text dates ("2023-05-11"), not real Excel date serials, so force
# text formatting before assignment to stop Excel's automatic date parsing.
$ws.Cells.Item($r, 25).NumberFormat = "@"
$ws.Cells.Item($r, 25).Value = "2023-05-11"                                  # Y  Startdatum
$ws.Cells.Item($r, 27).NumberFormat = "@"
$ws.Cells.Item($r, 27).Value = "2023-05-11"                                  # AA Slutdatum

$ws.Cells.Item($r, 29).Value = "miljön mölig mentroligen förändrad sedan förra observationen"  # AC Publik kommentar
$ws.Cells.Item($r, 30).Value = $true                                         # AD Ej återfunnen
$ws.Cells.Item($r, 31).Value = $false                                        # AE Osäker artbestämning
$ws.Cells.Item($r, 32).Value = ""                                            # AF Bestämningsmetod
$ws.Cells.Item($r, 33).Value = $false                                        # AG Ospontan

$ws.Cells.Item($r, 46).Value = ""                                            # AT Bestämningsår

$ws.Cells.Item($r, 49).Value = "Charlotte Wigermo"                                          # AW Rapportör
$ws.Cells.Item($r, 50).Value = "Ivar Anderberg, Bo Andersson, Barbro Ahlner"                 # AX Observatörer
$ws.Cells.Item($r, 51).Value = "Floraväkteri Sverige"                                        # AY Projektnamn
